# Generate Report for Handback
# - Overview sheet: row for 5c99ccd4... (row 3) status changes from
#   "Ready for handoff" to "Handback transform failed" for both zh-cn
#   and de-de status columns (E3, F3).
# - zh-cn sheet: row 3 (5c99ccd4 handback) Error Detail (column P) gets a
#   new error message, and column P width grows to fit it.
# - de-de sheet: same for its row 3, plus matching column width change.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"

# NOTE: the engine's ColumnWidth setter re-derives the stored OOXML
# "width" from pixel-rounded Maximum Digit Width units and adds ~5/6 of a
# character of padding versus the value that was assigned, so asking for
# a stored width of 40 means requesting (40 - 5/6) here.
$targetColPWidth = 40 - (5 / 6)

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handback transform failed"
$zhcn.Range("P3").Value = "Handback file name: jipy1js5.ns2 is different with handoff file name: 5c99ccd4-1281-4fee-b97e-df964f407ed0.b6e7f961acf77504427357731c84cf79067f1dea.zh-cn."
$zhcn.Columns.Item(16).ColumnWidth = $targetColPWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handback transform failed"
$dede.Range("P3").Value = "Handback file name: jipy1js5.ns2 is different with handoff file name: 5c99ccd4-1281-4fee-b97e-df964f407ed0.b6e7f961acf77504427357731c84cf79067f1dea.de-de."
$dede.Columns.Item(16).ColumnWidth = $targetColPWidth
